$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.38095238095238
$ws.Range("D3").Value = 1.616002906514979
$ws.Range("G3").Value = 31.34920634920635
$ws.Range("D4").Value = 1.616002906514979
$ws.Range("G4").Value = 31.34920634920635
$ws.Range("D5").Value = 1.616739367375187
$ws.Range("G5").Value = 31.31746031746032
$ws.Range("D6").Value = 1.616709792844247
$ws.Range("G6").Value = 31.3015873015873
$ws.Range("D7").Value = 1.616709792844247
$ws.Range("G7").Value = 31.3015873015873
$ws.Range("D8").Value = 1.616709792844247
$ws.Range("G8").Value = 31.3015873015873
$ws.Range("D9").Value = 1.616709792844247
$ws.Range("G9").Value = 31.3015873015873
$ws.Range("D10").Value = 1.616709792844247
$ws.Range("G10").Value = 31.3015873015873
$ws.Range("D11").Value = 1.617119590608179
$ws.Range("G11").Value = 31.28571428571428
$ws.Range("D12").Value = 1.617119590608179
$ws.Range("G12").Value = 31.28571428571428
$ws.Range("D13").Value = 1.617524586093871
$ws.Range("G13").Value = 31.26984126984127
$ws.Range("D14").Value = 1.618889793767077
$ws.Range("G14").Value = 31.20634920634921
$ws.Range("D15").Value = 1.617602953998483
$ws.Range("G15").Value = 31.17460317460317
$ws.Range("D16").Value = 1.616076088720251
$ws.Range("G16").Value = 31.12698412698413
$ws.Range("D17").Value = 1.617468969410897
$ws.Range("G17").Value = 31.06349206349206
$ws.Range("D18").Value = 1.618060995702578
$ws.Range("G18").Value = 31.03174603174603
$ws.Range("D19").Value = 1.614010785510716
$ws.Range("G19").Value = 30.93650793650794
$ws.Range("D20").Value = 1.614743065306736
$ws.Range("G20").Value = 30.88888888888889
$ws.Range("D21").Value = 1.615585135160913
$ws.Range("G21").Value = 30.84126984126984
$ws.Range("D22").Value = 1.614785057574994
$ws.Range("G22").Value = 30.82539682539683
$ws.Range("D23").Value = 1.615473819241175
$ws.Range("G23").Value = 30.76190476190476
$ws.Range("D24").Value = 1.610220546530715
$ws.Range("G24").Value = 30.6984126984127
$ws.Range("D25").Value = 1.603627636421485
$ws.Range("G25").Value = 30.55555555555556
$ws.Range("D26").Value = 1.599852271604728
$ws.Range("G26").Value = 30.44444444444444
$ws.Range("D27").Value = 1.599873141635445
$ws.Range("G27").Value = 30.28571428571428
$ws.Range("D28").Value = 1.600770944597433
$ws.Range("G28").Value = 30.20634920634921
$ws.Range("D29").Value = 1.5995685659886
$ws.Range("G29").Value = 30.09523809523809
$ws.Range("D30").Value = 1.60009321694631
$ws.Range("G30").Value = 29.80952380952381
$ws.Range("D31").Value = 1.598840566934998
$ws.Range("G31").Value = 29.53968253968254
$ws.Range("D32").Value = 1.598619831636245
$ws.Range("G32").Value = 29.34920634920635
$ws.Range("D33").Value = 1.599162124009779
$ws.Range("G33").Value = 29.11111111111111
$ws.Range("D34").Value = 1.599417458278982
$ws.Range("G34").Value = 28.98412698412698
$ws.Range("D35").Value = 1.597281007388254
$ws.Range("G35").Value = 28.79365079365079
$ws.Range("D36").Value = 1.603954294558374
$ws.Range("G36").Value = 28.28571428571428
$ws.Range("D37").Value = 1.606798331626519
$ws.Range("G37").Value = 27.98412698412698
$ws.Range("D38").Value = 1.610332034086212
$ws.Range("G38").Value = 27.50793650793651
$ws.Range("D39").Value = 1.605934149064359
$ws.Range("G39").Value = 27.15873015873016
$ws.Range("D40").Value = 1.603339620435709
$ws.Range("G40").Value = 26.6031746031746
$ws.Range("D41").Value = 1.607803722711402
$ws.Range("G41").Value = 25.6984126984127
$ws.Range("D42").Value = 1.603424777184309
$ws.Range("G42").Value = 24.57142857142857
$ws.Range("D43").Value = 1.560339712183514
$ws.Range("G43").Value = 23.36507936507936
$ws.Range("D44").Value = 1.501951814193386
$ws.Range("G44").Value = 22.09523809523809
$ws.Range("D45").Value = 1.428047081826624
$ws.Range("G45").Value = 19.6031746031746
$ws.Range("D46").Value = 1.315880376116133
$ws.Range("G46").Value = 16.17460317460317
$ws.Range("D47").Value = 1.215743644612701
$ws.Range("G47").Value = 11.65079365079365
$ws.Range("D48").Value = 1.137424104047591
$ws.Range("G48").Value = 7.158730158730159
$ws.Range("D49").Value = 1.148266439985296
$ws.Range("G49").Value = 3.07936507936508
$ws.Range("D50").Value = 1.401531578137174
$ws.Range("G50").Value = 0.7619047619047619

# D51 is cleared to an empty (inline string) cell, matching C51
$ws.Range("D51").Value = ""
